# Auto-generated edit script applying the diff to Sheets/Aegis_Profits.xlsx
# Updates cached numeric values (currentAveragePrice / LevePrice / LeveProfit columns)
# for specific Leve rows across the ALC, ARM, BSM, CRP, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# ALC row 6: Days of Chunder / Antidote
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 638361.75
$ws.Range("I6").Value = 1435150.2
$ws.Range("K6").Value = 4305450.6
$ws.Range("M6").Value = -4305338.6

# ALC row 9: Distill, My Heart / Distilled Water
$ws.Range("H9").Value = 277.14285
$ws.Range("I9").Value = 374
$ws.Range("K9").Value = 374
$ws.Range("M9").Value = -205

# ALC row 12: Don't Be So Tallow / Beeswax
$ws.Range("H12").Value = 62730
$ws.Range("I12").Value = 475
$ws.Range("J12").Value = 83481.664
$ws.Range("K12").Value = 475
$ws.Range("L12").Value = 83481.664
$ws.Range("M12").Value = -305
$ws.Range("N12").Value = -83821.664

# ALC row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 2193.2
$ws.Range("I62").Value = 2159.111
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 2159.111
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = -1535.111
$ws.Range("N62").Value = -3748

# ALC row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 2193.2
$ws.Range("I65").Value = 2159.111
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 10795.555
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = -7675.555
$ws.Range("N65").Value = -18740

# ALC row 88: The Grave of Hemlock Groves / Growth Formula Zeta
$ws.Range("H88").Value = 2105.087
$ws.Range("I88").Value = 769
$ws.Range("J88").Value = 2689.625
$ws.Range("K88").Value = 769
$ws.Range("L88").Value = 2689.625
$ws.Range("M88").Value = -363
$ws.Range("N88").Value = -3501.625

# ALC row 91: Dappling the Highlands (L) / Growth Formula Zeta
$ws.Range("H91").Value = 2105.087
$ws.Range("I91").Value = 769
$ws.Range("J91").Value = 2689.625
$ws.Range("K91").Value = 769
$ws.Range("L91").Value = 2689.625
$ws.Range("M91").Value = 635
$ws.Range("N91").Value = -5497.625

# ALC row 103: Let Loose the Juice / Persimmon Tannin
$ws.Range("H103").Value = 2779.375
$ws.Range("J103").Value = 2675.25
$ws.Range("L103").Value = 8025.75
$ws.Range("N103").Value = -9197.75

# ALC row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3154.5845
$ws.Range("I138").Value = 1860.4166
$ws.Range("J138").Value = 3740.6226
$ws.Range("K138").Value = 5581.2498
$ws.Range("L138").Value = 11221.8678
$ws.Range("M138").Value = -441.2497999999996
$ws.Range("N138").Value = -21501.8678

# ARM row 44: Very Slow Array / Mythril Plate
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 5820
$ws.Range("J44").Value = 5820
$ws.Range("L44").Value = 5820
$ws.Range("N44").Value = -6796

# ARM row 55: Employee Retention / Mythril Elmo
$ws.Range("H55").Value = 14645.714
$ws.Range("J55").Value = 15586.667
$ws.Range("L55").Value = 15586.667
$ws.Range("N55").Value = -16216.667

# BSM row 86: Through Thick and Thin / Adamantite Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 185626
$ws.Range("I86").Value = 185626
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 185626
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -184503
$ws.Range("N86").ClearContents()

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 185626
$ws.Range("I89").Value = 185626
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 928130
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -922514
$ws.Range("N89").ClearContents()

# CRP row 62: Splinter in the Sewers / Cedar Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2675
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# CRP row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 2675
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# GSM row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 125002536
$ws.Range("I97").Value = 142859900
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 142859900
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -142859404
$ws.Range("N97").Value = -1992

# LTW row 22: Skin off Their Backs / Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1239.8667
$ws.Range("I22").Value = 2675
$ws.Range("J22").Value = 718
$ws.Range("K22").Value = 2675
$ws.Range("L22").Value = 718
$ws.Range("M22").Value = -2380
$ws.Range("N22").Value = -1308

# LTW row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 1239.8667
$ws.Range("I27").Value = 2675
$ws.Range("J27").Value = 718
$ws.Range("K27").Value = 2675
$ws.Range("L27").Value = 718
$ws.Range("M27").Value = -2568
$ws.Range("N27").Value = -932

# LTW row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 69483.87
$ws.Range("I40").Value = 202600.6
$ws.Range("J40").Value = 2925.5
$ws.Range("K40").Value = 202600.6
$ws.Range("L40").Value = 2925.5
$ws.Range("M40").Value = -202464.6
$ws.Range("N40").Value = -3197.5

# LTW row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 3108.6667
$ws.Range("I68").Value = 2127.889
$ws.Range("J68").Value = 8993.333000000001
$ws.Range("K68").Value = 2127.889
$ws.Range("L68").Value = 8993.333000000001
$ws.Range("M68").Value = -1378.889
$ws.Range("N68").Value = -10491.333

# LTW row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 3108.6667
$ws.Range("I71").Value = 2127.889
$ws.Range("J71").Value = 8993.333000000001
$ws.Range("K71").Value = 10639.445
$ws.Range("L71").Value = 44966.665
$ws.Range("M71").Value = -6895.445
$ws.Range("N71").Value = -52454.665

# WVR row 15: Workplace Safety / Cotton Scarf
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9772.223
$ws.Range("J15").Value = 9817.647000000001
$ws.Range("L15").Value = 9817.647000000001
$ws.Range("N15").Value = -10393.647

# WVR row 54: No Country for Cold Men / Woolen Tights
$ws.Range("H54").Value = 7070
$ws.Range("I54").Value = 7070
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 7070
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -6550
$ws.Range("N54").ClearContents()

# WVR row 62: Pride Up in Smoke / Rainbow Cloth
$ws.Range("H62").Value = 6995233.5
$ws.Range("I62").Value = 12822345
$ws.Range("K62").Value = 12822345
$ws.Range("M62").Value = -12821721

# WVR row 65: Desperate for Diversionaries (L) / Rainbow Cloth
$ws.Range("H65").Value = 6995233.5
$ws.Range("I65").Value = 12822345
$ws.Range("K65").Value = 64111725
$ws.Range("M65").Value = -64108605

# WVR row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Range("H81").Value = 168937.42
$ws.Range("I81").Value = 126156.25
$ws.Range("J81").Value = 254499.75
$ws.Range("K81").Value = 252312.5
$ws.Range("L81").Value = 508999.5
$ws.Range("M81").Value = -251251.5
$ws.Range("N81").Value = -511121.5

# WVR row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Range("H84").Value = 168937.42
$ws.Range("I84").Value = 126156.25
$ws.Range("J84").Value = 254499.75
$ws.Range("K84").Value = 1261562.5
$ws.Range("L84").Value = 2544997.5
$ws.Range("M84").Value = -1256258.5
$ws.Range("N84").Value = -2555605.5
